$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.8668301361480815
$ws.Cells.Item(2, 3).Value = 0.3913830654369832
$ws.Cells.Item(2, 5).Value = 0.2911693146178109
$ws.Cells.Item(2, 6).Value = 1.693958096448057
$ws.Cells.Item(2, 7).Value = 0.4894475107197138
$ws.Cells.Item(2, 8).Value = 0.6351169897501734
$ws.Cells.Item(2, 10).Value = 0.03785625716331964
$ws.Cells.Item(2, 12).Value = 0.4733758092774849
$ws.Cells.Item(2, 13).Value = 0.2872598945466933
$ws.Cells.Item(2, 14).Value = 1.371259889002964
$ws.Cells.Item(2, 15).Value = 2.211005508782762
$ws.Cells.Item(3, 2).Value = 0.7913447093798709
$ws.Cells.Item(3, 3).Value = 0.3879792132131143
$ws.Cells.Item(3, 5).Value = 0.2926701382564314
$ws.Cells.Item(3, 6).Value = 1.695403417438818
$ws.Cells.Item(3, 7).Value = 0.4895540128132652
$ws.Cells.Item(3, 8).Value = 0.6388070456828814
$ws.Cells.Item(3, 10).Value = 0.03589546736919758
$ws.Cells.Item(3, 12).Value = 0.4682654592467514
$ws.Cells.Item(3, 13).Value = 0.2731674853945449
$ws.Cells.Item(3, 14).Value = 1.380507341196122
$ws.Cells.Item(3, 15).Value = 2.218697368680481
$ws.Cells.Item(4, 2).Value = 0.7450806083366217
$ws.Cells.Item(4, 3).Value = 0.385907596297983
$ws.Cells.Item(4, 5).Value = 0.293680579179485
$ws.Cells.Item(4, 6).Value = 1.697154269175819
$ws.Cells.Item(4, 7).Value = 0.4899658110189264
$ws.Cells.Item(4, 8).Value = 0.6413583086783348
$ws.Cells.Item(4, 10).Value = 0.03468315429760338
$ws.Cells.Item(4, 12).Value = 0.465288173938589
$ws.Cells.Item(4, 13).Value = 0.2645912384408859
$ws.Cells.Item(4, 14).Value = 1.386666153842171
$ws.Cells.Item(4, 15).Value = 2.224743154001686
$ws.Cells.Item(5, 2).Value = 0.7262502167925788
$ws.Cells.Item(5, 3).Value = 0.3850681315511366
$ws.Cells.Item(5, 5).Value = 0.2941147443753245
$ws.Cells.Item(5, 6).Value = 1.698085079305656
$ws.Cells.Item(5, 7).Value = 0.4902206704576386
$ws.Cells.Item(5, 8).Value = 0.6424698272742617
$ws.Cells.Item(5, 10).Value = 0.03418704770137637
$ws.Cells.Item(5, 12).Value = 0.4641154544250554
$ws.Cells.Item(5, 13).Value = 0.2611158948211525
$ws.Cells.Item(5, 14).Value = 1.389296985060895
$ws.Cells.Item(5, 15).Value = 2.227539509889056
$ws.Cells.Item(6, 2).Value = 0.7231248516869471
$ws.Cells.Item(6, 3).Value = 0.3849290286063365
$ws.Cells.Item(6, 5).Value = 0.2941881914067377
$ws.Cells.Item(6, 6).Value = 1.698252771620893
$ws.Cells.Item(6, 7).Value = 0.490268245471718
$ws.Cells.Item(6, 8).Value = 0.6426587357232023
$ws.Cells.Item(6, 10).Value = 0.03410454482709113
$ws.Cells.Item(6, 12).Value = 0.4639231799114185
$ws.Cells.Item(6, 13).Value = 0.2605400059929224
$ws.Cells.Item(6, 14).Value = 1.389741149207623
$ws.Cells.Item(6, 15).Value = 2.228023933342669
$ws.Cells.Item(7, 2).Value = 0.7448265616833396
$ws.Cells.Item(7, 3).Value = 0.3858962556483618
$ws.Cells.Item(7, 5).Value = 0.2936863437231931
$ws.Cells.Item(7, 6).Value = 1.697165942192058
$ws.Cells.Item(7, 7).Value = 0.4899688957646973
$ws.Cells.Item(7, 8).Value = 0.6413730079848534
$ws.Cells.Item(7, 10).Value = 0.03467647201143365
$ws.Cells.Item(7, 12).Value = 0.465272193798981
$ws.Cells.Item(7, 13).Value = 0.2645442892030303
$ws.Cells.Item(7, 14).Value = 1.386701143757982
$ws.Cells.Item(7, 15).Value = 2.224779519775893
$ws.Cells.Item(8, 2).Value = 0.8407861003877031
$ws.Cells.Item(8, 3).Value = 0.3902056760042996
$ws.Cells.Item(8, 5).Value = 0.2916683655288121
$ws.Cells.Item(8, 6).Value = 1.694277362816699
$ws.Cells.Item(8, 7).Value = 0.4894122997775625
$ws.Cells.Item(8, 8).Value = 0.6363300859648575
$ws.Cells.Item(8, 10).Value = 0.03718193247985013
$ws.Cells.Item(8, 12).Value = 0.4715805634760812
$ws.Cells.Item(8, 13).Value = 0.2823851178843952
$ws.Cells.Item(8, 14).Value = 1.37434870370484
$ws.Cells.Item(8, 15).Value = 2.21338306850933
$ws.Cells.Item(9, 2).Value = 1.02957798016115
$ws.Cells.Item(9, 3).Value = 0.3987973800178253
$ws.Cells.Item(9, 5).Value = 0.2884149804777021
$ws.Cells.Item(9, 6).Value = 1.69545630977953
$ws.Cells.Item(9, 7).Value = 0.4910727635191279
$ws.Cells.Item(9, 8).Value = 0.628704560970121
$ws.Cells.Item(9, 10).Value = 0.04202761650287101
$ws.Cells.Item(9, 12).Value = 0.4852172251377169
$ws.Cells.Item(9, 13).Value = 0.3179671020490034
$ws.Cells.Item(9, 14).Value = 1.35393427786795
$ws.Cells.Item(9, 15).Value = 2.201534717473663
$ws.Cells.Item(10, 2).Value = 1.168601037914016
$ws.Cells.Item(10, 3).Value = 0.4051899202546423
$ws.Cells.Item(10, 5).Value = 0.2864514975528394
$ws.Cells.Item(10, 6).Value = 1.700487196284158
$ws.Cells.Item(10, 7).Value = 0.493976120601161
$ws.Cells.Item(10, 8).Value = 0.6244795679008917
$ws.Cells.Item(10, 10).Value = 0.04554551687598973
$ws.Cells.Item(10, 12).Value = 0.4959992746114636
$ws.Cells.Item(10, 13).Value = 0.3444603146920215
$ws.Cells.Item(10, 14).Value = 1.341249106349977
$ws.Cells.Item(10, 15).Value = 2.199238237231754
$ws.Cells.Item(11, 2).Value = 1.231903822725656
$ws.Cells.Item(11, 3).Value = 0.408114307597188
$ws.Cells.Item(11, 5).Value = 0.2856504440029966
$ws.Cells.Item(11, 6).Value = 1.703678843426999
$ws.Cells.Item(11, 7).Value = 0.4956637650685849
$ws.Cells.Item(11, 8).Value = 0.6228561200700966
$ws.Cells.Item(11, 10).Value = 0.0471365232056371
$ws.Cells.Item(11, 12).Value = 0.501068298984606
$ws.Cells.Item(11, 13).Value = 0.3565865617552149
$ws.Cells.Item(11, 14).Value = 1.335978885319577
$ws.Cells.Item(11, 15).Value = 2.199586846457493
$ws.Cells.Item(12, 2).Value = 1.255882374892906
$ws.Cells.Item(12, 3).Value = 0.4092239428287172
$ws.Cells.Item(12, 5).Value = 0.2853603167826577
$ws.Cells.Item(12, 6).Value = 1.705017131970692
$ws.Cells.Item(12, 7).Value = 0.4963556732529497
$ws.Cells.Item(12, 8).Value = 0.6222842457513735
$ws.Cells.Item(12, 10).Value = 0.04773763276673293
$ws.Cells.Item(12, 12).Value = 0.5030112455388434
$ws.Cells.Item(12, 13).Value = 0.3611888808570072
$ws.Cells.Item(12, 14).Value = 1.334055007124697
$ws.Cells.Item(12, 15).Value = 2.199919297946053
$ws.Cells.Item(13, 2).Value = 1.250717870620463
$ws.Cells.Item(13, 3).Value = 0.4089848656645216
$ws.Cells.Item(13, 5).Value = 0.2854222137287543
$ws.Cells.Item(13, 6).Value = 1.704723142864424
$ws.Cells.Item(13, 7).Value = 0.4962043074628042
$ws.Cells.Item(13, 8).Value = 0.622405502156937
$ws.Cells.Item(13, 10).Value = 0.04760823451041318
$ws.Cells.Item(13, 12).Value = 0.5025917588925353
$ws.Cells.Item(13, 13).Value = 0.3601972328521299
$ws.Cells.Item(13, 14).Value = 1.33446615525942
$ws.Cells.Item(13, 15).Value = 2.199838782549335
$ws.Cells.Item(14, 2).Value = 1.233876417923511
$ws.Cells.Item(14, 3).Value = 0.4082055539030875
$ws.Cells.Item(14, 5).Value = 0.2856263104316064
$ws.Cells.Item(14, 6).Value = 1.703786347263645
$ws.Cells.Item(14, 7).Value = 0.4957196294762554
$ws.Cells.Item(14, 8).Value = 0.6228082122555492
$ws.Cells.Item(14, 10).Value = 0.04718600449201915
$ws.Cells.Item(14, 12).Value = 0.5012276781829428
$ws.Cells.Item(14, 13).Value = 0.356964991264519
$ws.Cells.Item(14, 14).Value = 1.335819167514678
$ws.Cells.Item(14, 15).Value = 2.199610179819047
$ws.Cells.Item(15, 2).Value = 1.223561428863491
$ws.Cells.Item(15, 3).Value = 0.4077284901241853
$ws.Cells.Item(15, 5).Value = 0.2857530453926849
$ws.Cells.Item(15, 6).Value = 1.703229415853144
$ws.Cells.Item(15, 7).Value = 0.4954296329504615
$ws.Cells.Item(15, 8).Value = 0.6230604682828016
$ws.Cells.Item(15, 10).Value = 0.04692719704540593
$ws.Cells.Item(15, 12).Value = 0.5003951831042031
$ws.Cells.Item(15, 13).Value = 0.3549864913010126
$ws.Cells.Item(15, 14).Value = 1.336657279204331
$ws.Cells.Item(15, 15).Value = 2.199496259897472
$ws.Cells.Item(16, 2).Value = 1.164465137907257
$ws.Cells.Item(16, 3).Value = 0.404999123879378
$ws.Cells.Item(16, 5).Value = 0.2865056993997221
$ws.Cells.Item(16, 6).Value = 1.700296771431852
$ws.Cells.Item(16, 7).Value = 0.4938732183358638
$ws.Cells.Item(16, 8).Value = 0.6245916679049799
$ws.Cells.Item(16, 10).Value = 0.04544135090640822
$ws.Cells.Item(16, 12).Value = 0.4956712895501028
$ws.Cells.Item(16, 13).Value = 0.3436693047862818
$ws.Cells.Item(16, 14).Value = 1.341603586374887
$ws.Cells.Item(16, 15).Value = 2.199243498875518
$ws.Cells.Item(17, 2).Value = 1.128225877647878
$ws.Cells.Item(17, 3).Value = 0.4033288556000798
$ws.Cells.Item(17, 5).Value = 0.2869910031250669
$ws.Cells.Item(17, 6).Value = 1.698728860086661
$ws.Cells.Item(17, 7).Value = 0.4930124322546305
$ws.Cells.Item(17, 8).Value = 0.625607441717591
$ws.Cells.Item(17, 10).Value = 0.04452742660575382
$ws.Cells.Item(17, 12).Value = 0.4928152533132675
$ws.Cells.Item(17, 13).Value = 0.3367453935228539
$ws.Cells.Item(17, 14).Value = 1.344766054536827
$ws.Cells.Item(17, 15).Value = 2.199445356043128
$ws.Cells.Item(18, 2).Value = 1.107387799881394
$ws.Cells.Item(18, 3).Value = 0.4023697113884452
$ws.Cells.Item(18, 5).Value = 0.2872788122265746
$ws.Cells.Item(18, 6).Value = 1.697912051155754
$ws.Cells.Item(18, 7).Value = 0.4925518601363734
$ws.Cells.Item(18, 8).Value = 0.6262197885275782
$ws.Cells.Item(18, 10).Value = 0.04400088768815635
$ws.Cells.Item(18, 12).Value = 0.491188005371427
$ws.Cells.Item(18, 13).Value = 0.3327699602993661
$ws.Cells.Item(18, 14).Value = 1.346632121842411
$ws.Cells.Item(18, 15).Value = 2.199692596885683
$ws.Cells.Item(19, 2).Value = 1.100333425490533
$ws.Cells.Item(19, 3).Value = 0.4020452316116518
$ws.Cells.Item(19, 5).Value = 0.2873777505191093
$ws.Cells.Item(19, 6).Value = 1.697650100062731
$ws.Cells.Item(19, 7).Value = 0.4924018464580087
$ws.Cells.Item(19, 8).Value = 0.6264319460693599
$ws.Cells.Item(19, 10).Value = 0.04382246155269343
$ws.Cells.Item(19, 12).Value = 0.4906397099228457
$ws.Cells.Item(19, 13).Value = 0.3314251618923549
$ws.Cells.Item(19, 14).Value = 1.347272032815688
$ws.Cells.Item(19, 15).Value = 2.199798828753984
$ws.Cells.Item(20, 2).Value = 1.132083021060055
$ws.Cells.Item(20, 3).Value = 0.4035064989950143
$ws.Cells.Item(20, 5).Value = 0.2869384441215672
$ws.Cells.Item(20, 6).Value = 1.698886969778073
$ws.Cells.Item(20, 7).Value = 0.493100490204526
$ws.Cells.Item(20, 8).Value = 0.6254964028244387
$ws.Cells.Item(20, 10).Value = 0.04462480607803343
$ws.Cells.Item(20, 12).Value = 0.4931176837359459
$ws.Cells.Item(20, 13).Value = 0.3374817317345418
$ws.Cells.Item(20, 14).Value = 1.344424530586764
$ws.Cells.Item(20, 15).Value = 2.19941029507612
$ws.Cells.Item(21, 2).Value = 1.238822978715689
$ws.Cells.Item(21, 3).Value = 0.4084343970312716
$ws.Cells.Item(21, 5).Value = 0.2855660039339973
$ws.Cells.Item(21, 6).Value = 1.704057988772277
$ws.Cells.Item(21, 7).Value = 0.4958605566615546
$ws.Cells.Item(21, 8).Value = 0.6226887628066038
$ws.Cells.Item(21, 10).Value = 0.04731006112128
$ws.Cells.Item(21, 12).Value = 0.5016277077914992
$ws.Cells.Item(21, 13).Value = 0.3579141000051607
$ws.Cells.Item(21, 14).Value = 1.335419806125806
$ws.Cells.Item(21, 15).Value = 2.199671885394224
$ws.Cells.Item(22, 2).Value = 1.308624741012295
$ws.Cells.Item(22, 3).Value = 0.4116680232773433
$ws.Cells.Item(22, 5).Value = 0.2847460418271872
$ws.Cells.Item(22, 6).Value = 1.708193344553408
$ws.Cells.Item(22, 7).Value = 0.4979724097556044
$ws.Cells.Item(22, 8).Value = 0.6211037871548086
$ws.Cells.Item(22, 10).Value = 0.04905703012134666
$ws.Cells.Item(22, 12).Value = 0.5073259014686613
$ws.Cells.Item(22, 13).Value = 0.3713281447895653
$ws.Cells.Item(22, 14).Value = 1.329953390764587
$ws.Cells.Item(22, 15).Value = 2.201011210582976
$ws.Cells.Item(23, 2).Value = 1.271366980399932
$ws.Cells.Item(23, 3).Value = 0.4099410308892573
$ws.Cells.Item(23, 5).Value = 0.2851766367380399
$ws.Cells.Item(23, 6).Value = 1.705917131683805
$ws.Cells.Item(23, 7).Value = 0.4968170683956998
$ws.Cells.Item(23, 8).Value = 0.6219268575888037
$ws.Cells.Item(23, 10).Value = 0.04812538254688548
$ws.Cells.Item(23, 12).Value = 0.5042722526471408
$ws.Cells.Item(23, 13).Value = 0.364163403586339
$ws.Cells.Item(23, 14).Value = 1.332832643254605
$ws.Cells.Item(23, 15).Value = 2.200189454213842
$ws.Cells.Item(24, 2).Value = 1.130339217756614
$ws.Cells.Item(24, 3).Value = 0.4034261829322361
$ws.Cells.Item(24, 5).Value = 0.2869621786287198
$ws.Cells.Item(24, 6).Value = 1.69881522483891
$ws.Cells.Item(24, 7).Value = 0.4930605723528458
$ws.Cells.Item(24, 8).Value = 0.6255465151500488
$ws.Cells.Item(24, 10).Value = 0.04458078427898471
$ws.Cells.Item(24, 12).Value = 0.4929809090598667
$ws.Cells.Item(24, 13).Value = 0.3371488169588375
$ws.Cells.Item(24, 14).Value = 1.344578784281588
$ws.Cells.Item(24, 15).Value = 2.199425737496227
$ws.Cells.Item(25, 2).Value = 0.978444838160101
$ws.Cells.Item(25, 3).Value = 0.3964585997531032
$ws.Cells.Item(25, 5).Value = 0.2892199936559692
$ws.Cells.Item(25, 6).Value = 1.694405479109179
$ws.Cells.Item(25, 7).Value = 0.4903283507331793
$ws.Cells.Item(25, 8).Value = 0.6305253882793806
$ws.Cells.Item(25, 10).Value = 0.04072406832683839
$ws.Cells.Item(25, 12).Value = 0.4813935029359442
$ws.Cells.Item(25, 13).Value = 0.3082786820899557
$ws.Cells.Item(25, 14).Value = 1.359050076991849
$ws.Cells.Item(25, 15).Value = 2.203615136609614
